$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old rows 13-23 down to 14-24)
$ws.Rows(13).Insert()

# The newly inserted row 13 should mirror the B/C formatting of row 14 (style 2/3),
# with no value in column A.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A13").Clear()

$ws.Range("B10").Value2 = '1) Gerais: - Levar os estudantes a compreenderem os mecanismos de obtenção e análise de daos de variáveis de processo ,identificando as causas especiais de variação ( instabilidade), e causas comuns, de natureza aleatória. 2) Específicos: - Ao final do curso os educandos devem:? Saber identificar causas especiais de variação? Saber analisar os resultados propondo a condições que levem o processo a uma menor variabilidadeSaber determinar a capacidade do processo e utilizar as cartas de controle'
$ws.Range("C10").Value2 = '1) Gerais: - Levar os estudantes a compreenderem os mecanismos de obtenção e análise de daos de variáveis de processo ,identificando as causas especiais de variação ( instabilidade), e causas comuns, de natureza aleatória. 2) Específicos: - Ao final do curso os educandos devem:? Saber identificar causas especiais de variação? Saber analisar os resultados propondo a condições que levem o processo a uma menor variabilidadeSaber determinar a capacidade do processo e utilizar as cartas de controle'
$ws.Range("B13").Value2 = '5840535 - Messias Borges Silva'
$ws.Range("C13").Value2 = '5840535 - Messias Borges Silva'
$ws.Range("B14").Value2 = 'IntroduçãoCartas de ControleAnálise da Cacidade de ProcessosCartas EspeciasCasos Práticos'
$ws.Range("C14").Value2 = 'IntroduçãoCartas de ControleAnálise da Cacidade de ProcessosCartas EspeciasCasos Práticos'
$ws.Range("B16").Value2 = 'I - DESCRITIVO:INTRODUÇÃO- A importância do CEP- Potencialidades- Natureza da variação - Causas especiais e causas comunsCARTAS DE CONTROLE- Cartas X, R- Carta X- Carta P- Carta nP- Carta C- Carta UCARTAS DE CONTROLE ESPECIAIS- Amplitude móvel- Soma acumulada (CUSUM)ANÁLISE DE CAPACIDADE DOS PROCESSOS- Indice Co- Indice CpK- Indice PPKCASOS PRÁTICOS- Utilização de situações reais vivenciados em ambiente indus-trial.'
$ws.Range("C16").Value2 = 'I - DESCRITIVO:INTRODUÇÃO- A importância do CEP- Potencialidades- Natureza da variação - Causas especiais e causas comunsCARTAS DE CONTROLE- Cartas X, R- Carta X- Carta P- Carta nP- Carta C- Carta UCARTAS DE CONTROLE ESPECIAIS- Amplitude móvel- Soma acumulada (CUSUM)ANÁLISE DE CAPACIDADE DOS PROCESSOS- Indice Co- Indice CpK- Indice PPKCASOS PRÁTICOS- Utilização de situações reais vivenciados em ambiente indus-trial.'
$ws.Range("B19").Value2 = 'duas provas escritas'
$ws.Range("C19").Value2 = 'duas provas escritas'
$ws.Range("B20").Value2 = 'serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso.A média da disciplina será a média aritmética das duas provas.'
$ws.Range("C20").Value2 = 'serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso.A média da disciplina será a média aritmética das duas provas.'
$ws.Range("B21").Value2 = 'uma prova escrita com conteúdo de todo o semestre'
$ws.Range("C21").Value2 = 'uma prova escrita com conteúdo de todo o semestre'
$ws.Range("B22").Value2 = '1 - Ferramentas Estatísticas Básicas p/ o Gerenciamento de Processos. Maria Cristina C. Werkena. Edit. FCO, 19962) Controle Estatístico da Qualidade, 4ª edição. Douglas C. Mont gomery, 2006. Edit. LTC3- Statisticial Quality Control, 5ª edição. Fugeno L. Grant and Richard S. Leavenworth MC Graw Hill, 1987'
$ws.Range("C22").Value2 = '1 - Ferramentas Estatísticas Básicas p/ o Gerenciamento de Processos. Maria Cristina C. Werkena. Edit. FCO, 19962) Controle Estatístico da Qualidade, 4ª edição. Douglas C. Mont gomery, 2006. Edit. LTC3- Statisticial Quality Control, 5ª edição. Fugeno L. Grant and Richard S. Leavenworth MC Graw Hill, 1987'
